# Generate Report for Handback
# Update the Correspond Handoff Datetime / Correspond Handback DateTime
# values for the zh-cn and de-de report sheets to reflect the re-run
# handback timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-20 00:51:18"
$wsZhCn.Range("H2").Value = "2016-03-20 00:51:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-20 00:51:22"
$wsDeDe.Range("H2").Value = "2016-03-20 00:51:39"
